$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old row 1 (Philip / 394728739813) contents; row numbering of
# everything below stays the same (no shifting of rows).
$ws.Range("A1:F1").ClearContents()

# --- Fall 2022 block ---
# Row 4 (was POLS 1101 / CPSC 3121 / CPSC 4148)
$ws.Range("A4").Value = "PSYC 1101"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("D4").Value = 3
$ws.Range("E4:F4").ClearContents()

# Row 5 (was ARTH 2125 / CPSC 3165 / CPSC 4155)
$ws.Range("A5").Value = "POLS 1101"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CPSC 4135"
$ws.Range("D5").Value = 3
$ws.Range("E5:F5").ClearContents()

# Row 6 (was DSCI 3111 / KINS 3258)
$ws.Range("A6").Value = "PSYC 1105"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "CPSC 4148"
$ws.Range("D6").Value = 3

# Row 7 (was CPSC 4000 / CPSC 3415)
$ws.Range("A7").Value = "DSCI 3111"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4155"
$ws.Range("D7").Value = 3

# Row 8 is new
$ws.Range("A8").Value = "CPSC 3121"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "CPSC 4157"
$ws.Range("D8").Value = 3

# Row 9 is new
$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0

# --- Fall 2023 block ---
# Row 13 (was CYBR 4125 / CPSC 4175 / CPSC 4176)
$ws.Range("A13").Value = "CPSC 4175"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3
$ws.Range("E13:F13").ClearContents()

# Rows 14 and 15 are removed entirely.
$ws.Range("A14:F15").ClearContents()
